$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel's date serial epoch (1899-12-30) so Value2 (raw serial) round-trips to a real date.
$epoch = Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0

$lastRow = $ws.UsedRange.Rows.Count
if ($lastRow -lt 2) { $lastRow = 113 }

# Column A held Excel date-serial numbers (formatted via a custom yyyy-mm-dd style).
# The dataloader now wants plain YYYYMMDD integers with no special number format,
# so convert each value and strip the date style back to the workbook default.
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $serial = $cell.Value2
    if ($null -ne $serial -and $serial -ne "") {
        $days = [double]$serial
        $date = $epoch.AddDays($days)
        $yyyymmdd = [int]$date.ToString("yyyyMMdd")
        $cell.Value = $yyyymmdd
    }
}

$ws.Range("A2:A$lastRow").Style = "Normal"
